$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.094883918762207
$ws.Range("B1").Value = 4.978244304656982
$ws.Range("C1").Value = 6.698893070220947
$ws.Range("D1").Value = 8.912601470947266
$ws.Range("E1").Value = 5.956182479858398
